$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10568.143
$ws.Range("I64").Value = 8994
$ws.Range("K64").Value = 8994
$ws.Range("M64").Value = -8746
$ws.Range("H67").Value = 10568.143
$ws.Range("I67").Value = 8994
$ws.Range("K67").Value = 8994
$ws.Range("M67").Value = -8136
$ws.Range("H74").Value = 15799.8
$ws.Range("I74").Value = 13333.333
$ws.Range("J74").Value = 19499.5
$ws.Range("K74").Value = 13333.333
$ws.Range("L74").Value = 19499.5
$ws.Range("M74").Value = -12397.333
$ws.Range("N74").Value = -21371.5
$ws.Range("H77").Value = 15799.8
$ws.Range("I77").Value = 13333.333
$ws.Range("J77").Value = 19499.5
$ws.Range("K77").Value = 66666.66500000001
$ws.Range("L77").Value = 97497.5
$ws.Range("M77").Value = -61986.66500000001
$ws.Range("N77").Value = -106857.5
$ws.Range("H98").Value = 410933.2
$ws.Range("I98").Value = 2722.5
$ws.Range("K98").Value = 2722.5
$ws.Range("M98").Value = -1224.5
$ws.Range("H107").Value = 343.25
$ws.Range("I107").Value = 379.875
$ws.Range("K107").Value = 379.875
$ws.Range("M107").Value = 1540.125
$ws.Range("H113").Value = 7406.1665
$ws.Range("I113").Value = 5316.6665
$ws.Range("J113").Value = 9495.666999999999
$ws.Range("K113").Value = 5316.6665
$ws.Range("L113").Value = 9495.666999999999
$ws.Range("M113").Value = -2062.6665
$ws.Range("N113").Value = -16003.667
$ws.Range("H122").Value = 410933.2
$ws.Range("I122").Value = 2722.5
$ws.Range("K122").Value = 8167.5
$ws.Range("M122").Value = -5717.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6308.385
$ws.Range("I61").Value = 5743.381
$ws.Range("J61").Value = 8681.4
$ws.Range("K61").Value = 5743.381
$ws.Range("L61").Value = 8681.4
$ws.Range("M61").Value = -5531.381
$ws.Range("N61").Value = -9105.4
$ws.Range("H63").Value = 6338.5884
$ws.Range("I63").Value = 3961.75
$ws.Range("K63").Value = 3961.75
$ws.Range("M63").Value = -3275.75
$ws.Range("H66").Value = 6338.5884
$ws.Range("I66").Value = 3961.75
$ws.Range("K66").Value = 19808.75
$ws.Range("M66").Value = -16376.75
$ws.Range("H97").Value = 1408.2273
$ws.Range("I97").Value = 846.3684
$ws.Range("J97").Value = 4966.6665
$ws.Range("K97").Value = 846.3684
$ws.Range("L97").Value = 4966.6665
$ws.Range("M97").Value = -350.3684
$ws.Range("N97").Value = -5958.6665
$ws.Range("H122").Value = 4555.3335
$ws.Range("I122").Value = 4322.625
$ws.Range("J122").Value = 4821.2856
$ws.Range("K122").Value = 12967.875
$ws.Range("L122").Value = 14463.8568
$ws.Range("M122").Value = -10517.875
$ws.Range("N122").Value = -19363.8568
$ws.Range("H132").Value = 3604.85
$ws.Range("I132").Value = 2587.9443
$ws.Range("J132").Value = 12757
$ws.Range("K132").Value = 7763.8329
$ws.Range("L132").Value = 38271
$ws.Range("M132").Value = -5233.8329
$ws.Range("N132").Value = -43331
$ws.Range("H136").Value = 6308.385
$ws.Range("I136").Value = 5743.381
$ws.Range("J136").Value = 8681.4
$ws.Range("K136").Value = 17230.143
$ws.Range("L136").Value = 26044.2
$ws.Range("M136").Value = -14680.143
$ws.Range("N136").Value = -31144.2

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 795.8333
$ws.Range("I22").Value = 797
$ws.Range("J22").Value = 793.5
$ws.Range("K22").Value = 797
$ws.Range("L22").Value = 793.5
$ws.Range("M22").Value = -624
$ws.Range("N22").Value = -1139.5
$ws.Range("H50").Value = 65000
$ws.Range("J50").Value = 65000
$ws.Range("L50").Value = 65000
$ws.Range("N50").Value = -66148
$ws.Range("H86").Value = 2565.5151
$ws.Range("I86").Value = 1912.52
$ws.Range("J86").Value = 4606.125
$ws.Range("K86").Value = 1912.52
$ws.Range("L86").Value = 4606.125
$ws.Range("M86").Value = -789.52
$ws.Range("N86").Value = -6852.125
$ws.Range("H89").Value = 2565.5151
$ws.Range("I89").Value = 1912.52
$ws.Range("J89").Value = 4606.125
$ws.Range("K89").Value = 9562.6
$ws.Range("L89").Value = 23030.625
$ws.Range("M89").Value = -3946.6
$ws.Range("N89").Value = -34262.625
$ws.Range("H94").Value = 4101.875
$ws.Range("I94").Value = 3884.1667
$ws.Range("J94").Value = 4755
$ws.Range("K94").Value = 3884.1667
$ws.Range("L94").Value = 4755
$ws.Range("M94").Value = -3433.1667
$ws.Range("N94").Value = -5657

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3253.9
$ws.Range("I58").Value = 1544.2667
$ws.Range("K58").Value = 1544.2667
$ws.Range("M58").Value = -1341.2667
$ws.Range("H94").Value = 2680.6428
$ws.Range("J94").Value = 3651.1428
$ws.Range("L94").Value = 3651.1428
$ws.Range("N94").Value = -4553.1428
$ws.Range("H99").Value = 2548.5
$ws.Range("I99").Value = 2548.5
$ws.Range("K99").Value = 2548.5
$ws.Range("M99").Value = -1050.5
$ws.Range("H107").Value = 3442.2222
$ws.Range("I107").Value = 1663.8334
$ws.Range("J107").Value = 6999
$ws.Range("K107").Value = 1663.8334
$ws.Range("L107").Value = 6999
$ws.Range("M107").Value = 256.1666
$ws.Range("N107").Value = -10839
$ws.Range("H122").Value = 6290.2856
$ws.Range("I122").Value = 1830.1428
$ws.Range("J122").Value = 10750.429
$ws.Range("K122").Value = 5490.428400000001
$ws.Range("L122").Value = 32251.287
$ws.Range("M122").Value = -3040.428400000001
$ws.Range("N122").Value = -37151.287
$ws.Range("H126").Value = 2548.5
$ws.Range("I126").Value = 2548.5
$ws.Range("K126").Value = 7645.5
$ws.Range("M126").Value = -5175.5
$ws.Range("H132").Value = 4291.8945
$ws.Range("I132").Value = 3384.1765
$ws.Range("J132").Value = 12007.5
$ws.Range("K132").Value = 10152.5295
$ws.Range("L132").Value = 36022.5
$ws.Range("M132").Value = -7622.529500000001
$ws.Range("N132").Value = -41082.5
$ws.Range("H136").Value = 3253.9
$ws.Range("I136").Value = 1544.2667
$ws.Range("K136").Value = 4632.800099999999
$ws.Range("M136").Value = -2082.800099999999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 100003760
$ws.Range("J75").Value = 5802.5
$ws.Range("L75").Value = 17407.5
$ws.Range("N75").Value = -19403.5
$ws.Range("H78").Value = 100003760
$ws.Range("J78").Value = 5802.5
$ws.Range("L78").Value = 52222.5
$ws.Range("N78").Value = -62206.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H97").Value = 1495.1875
$ws.Range("I97").Value = 1134.7037
$ws.Range("J97").Value = 3441.8
$ws.Range("K97").Value = 1134.7037
$ws.Range("L97").Value = 3441.8
$ws.Range("M97").Value = -638.7037
$ws.Range("N97").Value = -4433.8
$ws.Range("H122").Value = 8333.243
$ws.Range("I122").Value = 9656.421
$ws.Range("J122").Value = 6936.5557
$ws.Range("K122").Value = 28969.263
$ws.Range("L122").Value = 20809.6671
$ws.Range("M122").Value = -26519.263
$ws.Range("N122").Value = -25709.6671
$ws.Range("H126").Value = 3768.4546
$ws.Range("I126").Value = 2190.182
$ws.Range("K126").Value = 6570.545999999999
$ws.Range("M126").Value = -4100.545999999999
$ws.Range("H132").Value = 2503.4
$ws.Range("I132").Value = 1532.6666
$ws.Range("K132").Value = 4597.9998
$ws.Range("M132").Value = -2067.9998

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3347.348
$ws.Range("I22").Value = 1963.5714
$ws.Range("J22").Value = 5499.8887
$ws.Range("K22").Value = 1963.5714
$ws.Range("L22").Value = 5499.8887
$ws.Range("M22").Value = -1668.5714
$ws.Range("N22").Value = -6089.8887
$ws.Range("H27").Value = 3347.348
$ws.Range("I27").Value = 1963.5714
$ws.Range("J27").Value = 5499.8887
$ws.Range("K27").Value = 1963.5714
$ws.Range("L27").Value = 5499.8887
$ws.Range("M27").Value = -1856.5714
$ws.Range("N27").Value = -5713.8887
$ws.Range("H100").Value = 12681.158
$ws.Range("I100").Value = 7738.375
$ws.Range("K100").Value = 7738.375
$ws.Range("M100").Value = -7197.375
$ws.Range("H122").Value = 258959.81
$ws.Range("I122").Value = 315041.06
$ws.Range("J122").Value = 15941
$ws.Range("K122").Value = 945123.1799999999
$ws.Range("L122").Value = 47823
$ws.Range("M122").Value = -942673.1799999999
$ws.Range("N122").Value = -52723
$ws.Range("H132").Value = 8432.817999999999
$ws.Range("I132").Value = 6145.8335
$ws.Range("J132").Value = 11177.2
$ws.Range("K132").Value = 18437.5005
$ws.Range("L132").Value = 33531.60000000001
$ws.Range("M132").Value = -15907.5005
$ws.Range("N132").Value = -38591.60000000001
$ws.Range("H136").Value = 3402.4
$ws.Range("J136").Value = 5836.8184
$ws.Range("L136").Value = 17510.4552
$ws.Range("N136").Value = -22610.4552

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H122").Value = 3271.3572
$ws.Range("I122").Value = 1629.4
$ws.Range("K122").Value = 4888.200000000001
$ws.Range("M122").Value = -2438.200000000001
$ws.Range("H125").Value = 93600
$ws.Range("J125").Value = 93600
$ws.Range("L125").Value = 93600
$ws.Range("N125").Value = -103440
$ws.Range("H132").Value = 4062.2727
$ws.Range("I132").Value = 2035
$ws.Range("K132").Value = 6105
$ws.Range("M132").Value = -3575

Write-Output "applied edits"